# Generate Report for Handback
# Updates the "Correspond Handoff Datetime" (E) and "Correspond Handback DateTime" (H)
# columns for the f31b3b96... file row on both the zh-cn and de-de sheets, reflecting
# the freshly-generated handback report timestamps.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E3").Value = "2016-03-12 06:32:55"
$wsZhCn.Range("H3").Value = "2016-03-12 06:33:09"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E3").Value = "2016-03-12 06:32:58"
$wsDeDe.Range("H3").Value = "2016-03-12 06:33:14"
